# Update "合肥-漫展信息.xlsx": a new 漫展/event row (2024-07-06 合肥·次元日记动漫游戏嘉年华)
# is inserted into the two full-listing sheets ("展览" and "全部类型") at row 11,
# pushing subsequent rows down by one; the running index in column A for every
# pushed-down row is bumped by 1; and a handful of "想去人数" (F column) view
# counts are refreshed to newer scraped values across all three data sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$value) {
    # Force literal text even for strings that look like dates/numbers
    # (e.g. "2024-07-06"), matching the inlineStr cells already in the sheet,
    # without leaving a stray custom number-format style behind.
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.ClearFormats()
}

function Insert-EventRow($ws) {
    # Row 11 becomes the new "合肥·次元日记动漫游戏嘉年华" event; everything that
    # used to be at row 11+ shifts down to row 12+.
    $ws.Rows.Item(11).Insert()

    # Restore column-A's bordered/bold/centered style on the freshly inserted
    # (blank, unstyled) row by copying formats from the row right below it.
    $ws.Cells.Item(12, 1).Copy()
    $ws.Cells.Item(11, 1).PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item(11, 1).Value2 = 10
    Set-TextCell $ws.Cells.Item(11, 2) "2024-07-06"
    Set-TextCell $ws.Cells.Item(11, 3) "合肥·次元日记动漫游戏嘉年华"
    Set-TextCell $ws.Cells.Item(11, 4) "徽州大道5558号(徽州大道与紫云路交口) 合肥方圆荟(滨湖店)"
    Set-TextCell $ws.Cells.Item(11, 5) "2024.07.06 10:00-07.06 17:00"
    $ws.Cells.Item(11, 6).Value2 = 2
    $ws.Cells.Item(11, 7).Value2 = 45
    Set-TextCell $ws.Cells.Item(11, 8) "https://show.bilibili.com/platform/detail.html?id=87201"
    Set-TextCell $ws.Cells.Item(11, 9) "//i0.hdslb.com/bfs/openplatform/202406/BhvxoidA1717762410463.jpeg"
}

function Renumber-ShiftedRows($ws, [int]$lastRow) {
    # Column A holds a plain literal running index (0-based), not a formula;
    # bump every row that got pushed down by the insert above.
    for ($r = 12; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r, 1).Value2 + 1
    }
}

function Apply-FUpdates($ws, $updates) {
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value2 = $updates[$row]
    }
}

# ---- Sheet "展览" (overview/exhibitions) ----
$wsExpo = $wb.Worksheets.Item("展览")
Insert-EventRow $wsExpo
Renumber-ShiftedRows $wsExpo 22
Apply-FUpdates $wsExpo @{
    2  = 1059
    3  = 342
    4  = 1457
    5  = 8655
    8  = 636
    12 = 3507
    14 = 357
    16 = 1094
    18 = 1108
    20 = 187
    21 = 2231
    22 = 45
}

# ---- Sheet "演出" (performances) ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Cells.Item(2, 6).Value2 = 34

# ---- Sheet "全部类型" (all types) ----
$wsAll = $wb.Worksheets.Item("全部类型")
Insert-EventRow $wsAll
Renumber-ShiftedRows $wsAll 23
Apply-FUpdates $wsAll @{
    2  = 1059
    3  = 342
    4  = 1457
    5  = 8655
    8  = 636
    12 = 3507
    14 = 357
    16 = 1094
    18 = 1108
    20 = 187
    21 = 2231
    22 = 34
    23 = 45
}
